$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 511822.72
$ws.Range("I15").Value = 511822.72
$ws.Range("K15").Value = 1535468.16
$ws.Range("M15").Value = -1535299.16

$ws.Range("H98").Value = 1215
$ws.Range("I98").Value = 931.4194
$ws.Range("K98").Value = 931.4194
$ws.Range("M98").Value = 566.5806

$ws.Range("H122").Value = 1215
$ws.Range("I122").Value = 931.4194
$ws.Range("K122").Value = 2794.2582
$ws.Range("M122").Value = -344.2582000000002

$ws.Range("H132").Value = 2883.3394
$ws.Range("I132").Value = 2764.975
$ws.Range("J132").Value = 3179.25
$ws.Range("K132").Value = 8294.924999999999
$ws.Range("L132").Value = 9537.75
$ws.Range("M132").Value = -5764.924999999999
$ws.Range("N132").Value = -14597.75

$ws.Range("H137").Value = 26025.162
$ws.Range("I137").Value = 1274.8518
$ws.Range("K137").Value = 3824.5554
$ws.Range("M137").Value = -1274.5554

$ws.Range("H141").Value = 608.8570999999999
$ws.Range("I141").Value = 523.1111
$ws.Range("J141").Value = 1573.5
$ws.Range("K141").Value = 1569.3333
$ws.Range("L141").Value = 4720.5
$ws.Range("M141").Value = 3610.6667
$ws.Range("N141").Value = -15080.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10421592
$ws.Range("I32").Value = 11767566
$ws.Range("K32").Value = 11767566
$ws.Range("M32").Value = -11767279

$ws.Range("H61").Value = 1172.5616
$ws.Range("I61").Value = 965
$ws.Range("K61").Value = 965
$ws.Range("M61").Value = -753

$ws.Range("H74").Value = 3163.532
$ws.Range("I74").Value = 3139.3489
$ws.Range("K74").Value = 3139.3489
$ws.Range("M74").Value = -2265.3489

$ws.Range("H77").Value = 3163.532
$ws.Range("I77").Value = 3139.3489
$ws.Range("K77").Value = 15696.7445
$ws.Range("M77").Value = -11328.7445

$ws.Range("H123").Value = 22914.5
$ws.Range("J123").Value = 22914.5
$ws.Range("L123").Value = 22914.5
$ws.Range("N123").Value = -32714.5

$ws.Range("H132").Value = 1387.8518
$ws.Range("I132").Value = 1100.65
$ws.Range("J132").Value = 2208.4285
$ws.Range("K132").Value = 3301.95
$ws.Range("L132").Value = 6625.2855
$ws.Range("M132").Value = -771.9500000000003
$ws.Range("N132").Value = -11685.2855

$ws.Range("H136").Value = 1172.5616
$ws.Range("I136").Value = 965
$ws.Range("K136").Value = 2895
$ws.Range("M136").Value = -345

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1191
$ws.Range("I134").Value = 1056.68
$ws.Range("J134").Value = 1862.6
$ws.Range("K134").Value = 3170.04
$ws.Range("L134").Value = 5587.799999999999
$ws.Range("M134").Value = -635.04
$ws.Range("N134").Value = -10657.8

$ws.Range("H141").Value = 20000
$ws.Range("I141").Value = 20000
$ws.Range("K141").Value = 20000
$ws.Range("M141").Value = -14820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2938.5715
$ws.Range("I31").Value = 1632.2963
$ws.Range("J31").Value = 4541.727
$ws.Range("K31").Value = 1632.2963
$ws.Range("L31").Value = 4541.727
$ws.Range("M31").Value = -1337.2963
$ws.Range("N31").Value = -5131.727

$ws.Range("H34").Value = 2938.5715
$ws.Range("I34").Value = 1632.2963
$ws.Range("J34").Value = 4541.727
$ws.Range("K34").Value = 1632.2963
$ws.Range("L34").Value = 4541.727
$ws.Range("M34").Value = -1430.2963
$ws.Range("N34").Value = -4945.727

$ws.Range("H58").Value = 1059.7
$ws.Range("I58").Value = 659.17645
$ws.Range("K58").Value = 659.17645
$ws.Range("M58").Value = -456.17645

$ws.Range("H99").Value = 2666.0312
$ws.Range("I99").Value = 2616.611
$ws.Range("J99").Value = 2729.5715
$ws.Range("K99").Value = 2616.611
$ws.Range("L99").Value = 2729.5715
$ws.Range("M99").Value = -1118.611
$ws.Range("N99").Value = -5725.5715

$ws.Range("H122").Value = 1759.2941
$ws.Range("I122").Value = 1117.3334
$ws.Range("J122").Value = 3300
$ws.Range("K122").Value = 3352.0002
$ws.Range("L122").Value = 9900
$ws.Range("M122").Value = -902.0001999999999
$ws.Range("N122").Value = -14800

$ws.Range("H126").Value = 2666.0312
$ws.Range("I126").Value = 2616.611
$ws.Range("J126").Value = 2729.5715
$ws.Range("K126").Value = 7849.833
$ws.Range("L126").Value = 8188.7145
$ws.Range("M126").Value = -5379.833
$ws.Range("N126").Value = -13128.7145

$ws.Range("H136").Value = 1059.7
$ws.Range("I136").Value = 659.17645
$ws.Range("K136").Value = 1977.52935
$ws.Range("M136").Value = 572.4706499999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 133.33333
$ws.Range("I33").Value = 149
$ws.Range("J33").Value = 113.75
$ws.Range("K33").Value = 894
$ws.Range("L33").Value = 682.5
$ws.Range("M33").Value = -611
$ws.Range("N33").Value = -1248.5

$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7872
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27360
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2137
$ws.Range("I132").Value = 1935.7736
$ws.Range("J132").Value = 2803.5625
$ws.Range("K132").Value = 5807.3208
$ws.Range("L132").Value = 8410.6875
$ws.Range("M132").Value = -3277.3208
$ws.Range("N132").Value = -13470.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2325.276
$ws.Range("I7").Value = 2096.4348
$ws.Range("J7").Value = 3202.5
$ws.Range("K7").Value = 2096.4348
$ws.Range("L7").Value = 3202.5
$ws.Range("M7").Value = -1984.4348
$ws.Range("N7").Value = -3426.5

$ws.Range("H40").Value = 2801
$ws.Range("I40").Value = 2135.4375
$ws.Range("J40").Value = 5843.5713
$ws.Range("K40").Value = 2135.4375
$ws.Range("L40").Value = 5843.5713
$ws.Range("M40").Value = -1999.4375
$ws.Range("N40").Value = -6115.5713

$ws.Range("H126").Value = 2325.276
$ws.Range("I126").Value = 2096.4348
$ws.Range("J126").Value = 3202.5
$ws.Range("K126").Value = 6289.3044
$ws.Range("L126").Value = 9607.5
$ws.Range("M126").Value = -3819.3044
$ws.Range("N126").Value = -14547.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1535.6471
$ws.Range("I132").Value = 987.7838
$ws.Range("J132").Value = 2983.5715
$ws.Range("K132").Value = 2963.3514
$ws.Range("L132").Value = 8950.7145
$ws.Range("M132").Value = -433.3514
$ws.Range("N132").Value = -14010.7145

$ws.Range("H136").Value = 556.73
$ws.Range("I136").Value = 438.0625
$ws.Range("J136").Value = 1031.4
$ws.Range("K136").Value = 1314.1875
$ws.Range("L136").Value = 3094.2
$ws.Range("M136").Value = 1235.8125
$ws.Range("N136").Value = -8194.200000000001
